$d = $word.ActiveDocument

# Locate the paragraph whose text is "Azure tenant with Admin Privileges"
# (the first bullet under the "Requirements" heading, list numId 6, level 0).
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a", "`n") -eq "Azure tenant with Admin Privileges") {
        $para = $p
        break
    }
}

# Insert a brand-new paragraph right after it, inheriting the same list
# (numId 6) / style / run formatting, then demote it to sub-bullet level 1
# and give it its own text.
$null = $para.Range.InsertParagraphAfter()

$newPara = $para.Next()
$newPara.Range.ListFormat.ListLevelNumber = 2
$newPara.Range.Text = "In case you are going to deploy to an already existing RG. The user running the Wizard needs to be owner of the Resource Group"
